# Sprint 2 / Entregaveis 2 - Planilha de riscos 1.xlsx
# Update the "Matriz de risco" header label and rebuild the risk-response
# table (rows 10-16) to add a "Plano de resposta" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Probability/Impact matrix header: rename the top score label ------
$ws.Range("F2").Value = "Muito Grave (10)"

# --- 2. Clear the old risk table (rows 10-16, cols C:E) before rebuilding -
$ws.Range("C10:E16").Clear()

# --- 3. Rebuild the header row (row 10) across B:E -------------------------
$ws.Range("B10").Value = "Riscos do projeto"
$ws.Range("C10").Value = "Pontuação"
$ws.Range("D10").Value = "Impacto"
$ws.Range("E10").Value = "Plano de resposta"

$headerRange = $ws.Range("B10:D10")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 12
$headerRange.HorizontalAlignment = -4108

$planoHeader = $ws.Range("E10")
$planoHeader.Font.Bold = $true
$planoHeader.Font.Size = 11
$planoHeader.HorizontalAlignment = -4108
$planoHeader.VerticalAlignment = -4108

# --- 4. Risk rows 11-16 -----------------------------------------------------
# columns: B = risk description, C = score, D = impact-color cell, E = response plan
$rows = @(
    @{ Row = 11; Risk = "Integrante sair da equipe"; Score = 10; Color = 65535;   Plan = "Reorganização de tarefas entre os remanescentes " },
    @{ Row = 12; Risk = "Dificuldade de compreender a lógica"; Score = 12; Color = 15773696; Plan = "Explicar novamente o motivo da dúvida" },
    @{ Row = 13; Risk = "Erro de código"; Score = 18; Color = 65535;   Plan = "Revisão do Código" },
    @{ Row = 14; Risk = "Falta de comunicação"; Score = 21; Color = 5287936;  Plan = "Aumentar comunicação entre os integrantes" },
    @{ Row = 15; Risk = "Integrante se recusar a colaborar com a equipe"; Score = 20; Color = 255; Plan = "Conversar com superiores para ajudar a resolver" },
    @{ Row = 16; Risk = "Entregas atrasadas"; Score = 27; Color = 255; Plan = "Cobrança maior sobre o atraso e cuidado para a não reincidencia" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $riskCell = $ws.Cells.Item($rowNum, 2)
    $riskCell.Value = $r.Risk
    $riskCell.HorizontalAlignment = -4108

    $scoreCell = $ws.Cells.Item($rowNum, 3)
    $scoreCell.Value = $r.Score
    $scoreCell.HorizontalAlignment = -4108

    $colorCell = $ws.Cells.Item($rowNum, 4)
    $colorCell.Interior.Color = $r.Color

    $planCell = $ws.Cells.Item($rowNum, 5)
    $planCell.Value = $r.Plan
    $planCell.HorizontalAlignment = -4108
}

# D13 also keeps a centered horizontal alignment (matches source formatting)
$ws.Range("D13").HorizontalAlignment = -4108
# E12 is centered vertically as well as horizontally
$ws.Range("E12").VerticalAlignment = -4108

# --- 5. Borders: every cell in B10:E16 gets an individual medium box -------
for ($rowNum = 10; $rowNum -le 16; $rowNum++) {
    for ($col = 2; $col -le 5; $col++) {
        $cell = $ws.Cells.Item($rowNum, $col)
        $cell.Borders.Item(7).LineStyle = 1
        $cell.Borders.Item(7).Weight = -4138
        $cell.Borders.Item(8).LineStyle = 1
        $cell.Borders.Item(8).Weight = -4138
        $cell.Borders.Item(9).LineStyle = 1
        $cell.Borders.Item(9).Weight = -4138
        $cell.Borders.Item(10).LineStyle = 1
        $cell.Borders.Item(10).Weight = -4138
    }
}

# --- 6. Row heights: rows 11-15 now match the thick-bottom style of 10/16 --
$ws.Rows.Item(10).RowHeight = 16
for ($rowNum = 11; $rowNum -le 16; $rowNum++) {
    $ws.Rows.Item($rowNum).RowHeight = 15
}

# --- 7. Clear the leftover formatting below the table (rows 17-20, C:D) ----
$ws.Range("C17:D20").Clear()

# --- 8. Column widths -------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.54296875
$ws.Columns.Item(2).ColumnWidth = 42.26953125
$ws.Columns.Item(3).ColumnWidth = 27.36328125
$ws.Columns.Item(5).ColumnWidth = 55.36328125
$ws.Columns.Item(6).ColumnWidth = 42.26953125

# --- 9. View state: zoom + selection ---------------------------------------
$ws.Range("E16").Select()
$excel.ActiveWindow.Zoom = 70

Write-Output "edit complete"
